$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" - refresh the handoff/handback datetime
# stamps for the first data row (514bbbb0-...) on the localized report
# sheets, simulating a newly generated handback report.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-23 22:45:32"
$wsZhCn.Range("K2").Value = "2016-08-23 22:45:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-23 22:45:37"
$wsDeDe.Range("K2").Value = "2016-08-23 22:45:55"
